# Localization: Content - Disabled all languages except english... AGAIN -____-
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tech")

# Disable the "android" (F) and "iOS" (G) flags for every language row
# except the first one (row 5, english). Rows 6-16 cover: french, italian,
# german, spanish, brazilian, russian, chinese, japanese, korean,
# traditional chinese and turkish.
for ($r = 6; $r -le 16; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
}

# Move the active selection on the sheet (as recorded in the saved view)
$ws.Range("H16").Select()
